# The diff updates the raw "start time" / "end time" timestamp values
# (stored as plain numbers, e.g. .NET ticks) on Sheet1. The dependent
# formulas in B3 (=B2-B1) and B4 (=B3/POWER(10,9)) recalculate
# automatically, picking up the new results (233988000000 / 233.988).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = [double]"1.5018677321969999E+18"
$ws.Range("B2").Value = [double]"1.5018679661849999E+18"
